$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()

# Update the "Sprint 4 (GDP-2)" section durations (row 34-38, column D)
$ws.Range("D34").Value = 100
$ws.Range("D35").Value = 140
$ws.Range("D36").Value = 140
$ws.Range("D37").Value = 115
$ws.Range("D38").Value = 105

# Move the selection/active cell to D39 (matches the post-edit selection)
$ws.Range("D39").Select()
